# Apply the "canine breed xl files with stat bar query" edit:
# Insert a new column B ("StatQuery" header + stat-bar Neo4j query), pushing the
# existing "dbExcel"/"WebExcel" columns one slot to the right, and move the
# active selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column B (dbExcel / file name column).
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "StatQuery"

# The new "stat bar" Neo4j query text, formatted like the existing query cell (A2).
$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['English Setter']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").Value = $statQuery

# Match the wrapped-text formatting used by the existing query cell (A2).
$ws.Range("B2").WrapText = $true

# Give the new column the same width as column A.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()

# Move the active selection to A2 (was B2 before the edit).
[void]$ws.Range("A2").Select()
